$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.146.96"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.222.77"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "2.566.10"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "2.219.19"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.736"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "40.068.02"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.96%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E39").Value = "  +4.70%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "2.074.60"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("E44").Value = "  +12.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0271"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("E47").Value = "  +4.44%  "
$ws.Range("E48").Value = "  -11.15%  "
$ws.Range("D49").Value = "2.438.12"
$ws.Range("E50").Value = "  +4.93%  "
$ws.Range("E51").Value = "  +1.70%  "
